# Updated run for publication
# Replaces the frequency table values in rows 2-5 (columns B:X) with the
# newly computed (unrounded) proportions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.00643440050219711,0.670590081607031,0.354048964218456,0.458254865034526,0.288292529817954,0.952762084118016,0.227244193345888,0.038135593220339,0.189579409918393,0.0309165097300691,0.496704331450094,0.996390458254865,0.200878844946642,0.0232266164469554,0.170590081607031,0.00125549278091651,0.207470182046453,0.00204017576898933,0.999529190207156,0.00156936597614564,0.168392969240427,0.0306026365348399,0.32391713747646)
    3 = @(0.0530445699937225,0.0164783427495292,0.174042686754551,0.00768989328311362,0.0205586942875078,0.00800376647834275,0.0091023226616447,0.00156936597614564,0.00109855618330195,0.290646578782172,0.179535467671061,0.000941619585687382,0.713747645951036,0.975360954174514,0.000156936597614564,0,0.495919648462021,0.997802887633396,0.000313873195229127,0.491368487131199,0.506591337099812,0.18361581920904,0.00109855618330195)
    4 = @(0.462021343377276,0.308066541117389,0.46045197740113,0.0397049591964846,0.519774011299435,0.0182046453232894,0.279190207156309,0.275266792215945,0.0186754551161331,0.168549905838041,0.0017263025737602,0.000627746390458255,0.000470809792843691,0.00125549278091651,0.0186754551161331,0.998744507219083,0.27510985561833,0,0.000156936597614564,0.0505335844318895,0.313088512241055,0.29331450094162,0.657407407407407)
    5 = @(0.471437539234149,0.00266792215944758,0.00784682988072819,0.491368487131199,0.170747018204645,0.019774011299435,0.481167608286252,0.684714375392341,0.790646578782172,0.505963590709353,0.321563088512241,0.00204017576898933,0.0773697426239799,0,0.810577526679222,0,0.0211864406779661,0.000156936597614564,0,0.456214689265537,0.0112994350282486,0.485247959824231,0.00800376647834275)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    $rng = $ws.Range("B${rowNum}:X${rowNum}")
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $rng.Value = $arr
}
